$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "List of parts updated": the goldpin / JST-XH connector rows that were
# previously marked "Ordered" (STATUS column C, rows 35-40) have now
# arrived, so flip their STATUS to "Ready" - matching the rest of the list.
$ws.Range("C35:C40").Value = "Ready"
